$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fitness (column C) values for rows 2-33 to 7295
$ws.Range("C2:C33").Value = 7295

# Update Fitness (column C) values for rows 34-116 to 7293
$ws.Range("C34:C116").Value = 7293
